$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the report title (A5) from "Daftar..." to "Laporan..."
$ws.Range("A5").Value = "Laporan Surat Izin Usaha Jasa Konstruksi (SIUJK) Bulan 01 Tahun 2015"

# 2. Insert six new rows after the last existing data row (row 13), pushing the
#    footer block down by 6 rows (16->22, 17->23, 21->27, 22->28).
for ($i = 0; $i -lt 6; $i++) {
  $ws.Rows.Item(14).Insert()
}

# Copy the formatting of the last data row (13) down across the 6 new rows so
# they look like the rest of the table (border, wrap, vertical centering).
$ws.Range("A13:L13").Copy()
$ws.Range("A14:A19").PasteSpecial(-4122)

# The "Ket" columns (J:L) need centered + wrapped text (style used by K9 etc.)
$jkl = $ws.Range("J14:L19")
$jkl.WrapText = $true
$jkl.HorizontalAlignment = -4108
$jkl.VerticalAlignment = -4108
$jkl.Borders.LineStyle = 1

# Match the row height used by the other data rows
$ws.Range("14:19").RowHeight = 30

# 3. Fill in the six new SIUJK records (rows 6-11 of the table)
$ws.Cells.Item(14,1).Value = 6
$ws.Cells.Item(14,2).Value = "1.005123.1110.2.00006"
$ws.Cells.Item(14,3).Value = "EFNI ZAHARA"
$ws.Cells.Item(14,4).Value = ' "CV. AWE BEUTARI"'
$ws.Cells.Item(14,5).Value = " Gp. Lhok Awe-awe"
$ws.Cells.Item(14,6).Value = "BANGUNAN GEDUNG, BANGUNAN SIPIL"
$ws.Cells.Item(14,7).Value = "KUALA"
$ws.Cells.Item(14,8).Value = "14 Januari 2015"
$ws.Cells.Item(14,9).Value = "03 November 2017"
$ws.Cells.Item(14,11).Value = "√"

$ws.Cells.Item(15,1).Value = 7
$ws.Cells.Item(15,2).Value = "1.012707.1110.2.00007"
$ws.Cells.Item(15,3).Value = "SURYADI NURDIN"
$ws.Cells.Item(15,4).Value = ' "CV.  SAUDARA SEPAKAT"'
$ws.Cells.Item(15,5).Value = "Jl. B. Aceh - Medan Gp. Matang Glp. Dua Mns. Dayah"
$ws.Cells.Item(15,6).Value = "BANGUNAN GEDUNG, BANGUNAN SIPIL"
$ws.Cells.Item(15,7).Value = "PEUSANGAN"
$ws.Cells.Item(15,8).Value = "15 Januari 2015"
$ws.Cells.Item(15,9).Value = "22 Mei 2018"
$ws.Cells.Item(15,12).Value = "√"

$ws.Cells.Item(16,1).Value = 8
$ws.Cells.Item(16,2).Value = "1.018010.1110.2.00008"
$ws.Cells.Item(16,3).Value = "ABDUL MANAF ISDA"
$ws.Cells.Item(16,4).Value = ' "CV. AMICO GLOBAL"'
$ws.Cells.Item(16,5).Value = "Dsn. Balee Aron Gp. Bireuen Meunasah Blang"
$ws.Cells.Item(16,6).Value = "BANGUNAN GEDUNG, BANGUNAN SIPIL"
$ws.Cells.Item(16,7).Value = "KOTA JUANG"
$ws.Cells.Item(16,8).Value = "15 Januari 2015"
$ws.Cells.Item(16,9).Value = "15 Januari 2018"
$ws.Cells.Item(16,10).Value = "√"

$ws.Cells.Item(17,1).Value = 9
$ws.Cells.Item(17,2).Value = "1.005255.1110.2.00009"
$ws.Cells.Item(17,3).Value = "H. SYAHRIL MUSHERI"
$ws.Cells.Item(17,4).Value = ' "CV. GAPU GAGAH GETNA"'
$ws.Cells.Item(17,5).Value = "Jl. Medan - B. Aceh Km. 220 Gp. Geulanggang Baro"
$ws.Cells.Item(17,6).Value = "BANGUNAN GEDUNG, JASA PELAKSANAAN SPESIALIS"
$ws.Cells.Item(17,7).Value = "KOTA JUANG"
$ws.Cells.Item(17,8).Value = "19 Januari 2015"
$ws.Cells.Item(17,9).Value = "30 Mei 2017"
$ws.Cells.Item(17,12).Value = "√"

$ws.Cells.Item(18,1).Value = 10
$ws.Cells.Item(18,2).Value = "1.007333.1110.2.00010"
$ws.Cells.Item(18,3).Value = "SYAHRUL RAMADHAN"
$ws.Cells.Item(18,4).Value = ' "CV. UCOT COMMUNITY"'
$ws.Cells.Item(18,5).Value = "Jl. Medan - B. Aceh Km. 220 Gp. Geulanggang Baro"
$ws.Cells.Item(18,6).Value = "BANGUNAN GEDUNG, BANGUNAN SIPIL"
$ws.Cells.Item(18,7).Value = "KOTA JUANG"
$ws.Cells.Item(18,8).Value = "22 Januari 2015"
$ws.Cells.Item(18,9).Value = "29 Mei 2017"
$ws.Cells.Item(18,11).Value = "√"

$ws.Cells.Item(19,1).Value = 11
$ws.Cells.Item(19,2).Value = "1.015876.1110.2.00011"
$ws.Cells.Item(19,3).Value = "MUNAWIR"
$ws.Cells.Item(19,4).Value = ' "CV. DAKOTA PERKASA"'
$ws.Cells.Item(19,5).Value = "Jl. Raya P. Ara Gp. Pulo Ara Geudong Teungoh"
$ws.Cells.Item(19,6).Value = "BANGUNAN GEDUNG, BANGUNAN SIPIL"
$ws.Cells.Item(19,7).Value = "KOTA JUANG"
$ws.Cells.Item(19,8).Value = "27 Januari 2015"
$ws.Cells.Item(19,9).Value = "19 Mei 2017"
$ws.Cells.Item(19,11).Value = "√"

# 4. Update the signature block text (now shifted to rows 22/23 and 27-29)
$ws.Cells.Item(22,8).Value = "Kepala Kantor Pelayanan Perizinan Terpadu Satu Pintu"
$ws.Cells.Item(23,8).Value = "Kabupaten Bireuen"
$ws.Cells.Item(27,8).Value = "MUHAMMAD NASIR, SP"
$ws.Cells.Item(28,8).Value = "PEMBINA"
$ws.Cells.Item(29,8).Value = "NIP. 19621231 198711 1 002"

# 5. Keep the active selection in sync with the new last row (matches the diff)
$ws.Range("L19").Select()
